$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) values for columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values for columns B:E
$ws.Range("B2").Value = 10.725280943399168
$ws.Range("C2").Value = 9.011600379844781
$ws.Range("D2").Value = 8.4015339261526361
$ws.Range("E2").Value = 0.031016668627394068

# Update row 3 values for columns B:E
$ws.Range("B3").Value = 28.315069971155811
$ws.Range("C3").Value = 3.2432977487955794
$ws.Range("D3").Value = 1.3557792669367603
$ws.Range("E3").Value = 0.95061384647132041

# Update selection to match new sqref B1:E3
$ws.Range("B1:E3").Select()
